$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") holds the value 49 for every data row (rows 2-51).
# Recover the dropped data by restoring the correct competition ID: 249.
$ws.Range("B2:B51").Value = 249
